$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.343.04"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "3.887.29"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("D7").Value = "3.887.33"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").Value = "4.539.35"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "3.874.37"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "69.435.68"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +10.03%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000165"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.16%  "
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "4.035.94"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "3.853.20"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.97%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.54%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.325"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "434.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +22.70%  "
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.99%  "
